$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "Dia 1: ..." paragraph: merge the split runs (and drop the proofErr
#    spell/grammar markers) back into a single run.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(2)
$xml = '<w:p ' + $wNs + '><w:r><w:t>Dia 1: Establecer las bases y funcionalidades de la app.</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) "Dia 2: ..." paragraph: reword into three runs, then add a brand new
#    "Dia 3: ..." paragraph (underlined paragraph mark) right after it.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(3)
$xml = '<w:p ' + $wNs + '>' +
       '<w:r><w:t xml:space="preserve">Dia 2: </w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">Instalación del entorno. </w:t></w:r>' +
       '<w:r><w:t>Subir el proyecto a GitHub. Definir modelo de datos.</w:t></w:r>' +
       '</w:p>'
$p.Range.InsertXML($xml)

$p = $d.Paragraphs(3)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$newPara = $d.Paragraphs(4)
$xml = '<w:p ' + $wNs + '>' +
       '<w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
       '<w:r><w:t>Dia 3: Instalación de Angular material.</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> Creación módulo login.</w:t></w:r>' +
       '</w:p>'
$newPara.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3) "- Login y registro ..." paragraph: merge the split runs (drop the
#    spell-check proofErr markers around "Login").
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(7)
$xml = '<w:p ' + $wNs + '><w:r><w:t>- Login y registro para usuarios y administradores.</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4) "Administrador: boolean" paragraph: keep the "Administrador" run, merge
#    the remaining two runs (drop the proofErr markers around "boolean").
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(24)
$xml = '<w:p ' + $wNs + '>' +
       '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
       '<w:r><w:t>Administrador</w:t></w:r>' +
       '<w:r><w:t>: boolean</w:t></w:r>' +
       '</w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 5) "Nombre" paragraph (under "Hora fin"): gains a lastRenderedPageBreak.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(32)
$xml = '<w:p ' + $wNs + '>' +
       '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
       '<w:r><w:lastRenderedPageBreak/><w:t>Nombre</w:t></w:r>' +
       '</w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 6) "Disponible: boolean" paragraph: merge the split runs (drop the
#    proofErr markers around "boolean").
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(33)
$xml = '<w:p ' + $wNs + '>' +
       '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
       '<w:r><w:t>Disponible: boolean</w:t></w:r>' +
       '</w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 7) "Capacidad" paragraph: loses its lastRenderedPageBreak (it moved up to
#    the "Nombre" paragraph above, step 5).
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(34)
$xml = '<w:p ' + $wNs + '>' +
       '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
       '<w:r><w:t>Capacidad</w:t></w:r>' +
       '</w:p>'
$p.Range.InsertXML($xml)

Write-Output "Done"
